# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G (header "K") values change for rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 6
    6  = 3
    7  = 10
    8  = 3
    9  = 3
    10 = 8
    11 = 8
    12 = 4
    13 = 6
    14 = 7
    15 = 5
    16 = 4
    17 = 5
    18 = 3
    19 = 8
    20 = 6
    21 = 4
    22 = 5
    23 = 10
    24 = 3
    25 = 3
    26 = 8
    27 = 9
    28 = 4
    29 = 6
    30 = 3
    31 = 4
    32 = 1
    33 = 4
    34 = 2
    35 = 2
    36 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
